# moved tbot to slash commands
# Adds new "custom status" rows (291-298) sourced from the Tbot -> slash-command
# migration, widens column B slightly, and keeps the sheet dimension in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 291; Status = 'Spinning the wheel for which widely-considered-mediocre card is now busted because of a certain tech'; Creator = 'stingray' }
    @{ Row = 292; Status = 'Getting mass downvoted after explaining why ANB isn''t S tier'; Creator = '4-Pound Plate of Pasta' }
    @{ Row = 293; Status = 'Gee i sure hope I don''t miss this eyespore fusion'; Creator = '4-Pound Plate of Pasta' }
    @{ Row = 294; Status = 'Reading the 400th crippling DMD nerf on the sub'; Creator = '4-Pound Plate of Pasta' }
    @{ Row = 295; Status = 'Roping my opponent so I can listen to the low health theme'; Creator = '4-Pound Plate of Pasta' }
    @{ Row = 296; Status = 'Running Buff-Shroom in Cycle-Cap to counter ZM Sig'; Creator = '4-Pound Plate of Pasta' }
    @{ Row = 297; Status = 'Getting a Magic Beanstalk from the 15 daily streak reward'; Creator = '4-Pound Plate of Pasta' }
    @{ Row = 298; Status = 'getting nothing because it''s not worth logging in'; Creator = 'Stingray' }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Status
    $ws.Cells.Item($r.Row, 2).Value = $r.Creator
}

# Column B grows from 20.83 to 22.83 characters wide to fit the new entries.
$ws.Columns.Item(2).ColumnWidth = 22

